# Write transient segment data for .sfr file (ITEM6d sheet).
#
# This renumbers the existing "segment" blocks (column A) from 1,3,4,5,6,7
# down to 1,2,3,4,5,6 (each segment has an "x" row and a "z" row), then
# appends a new segment-7 x/z block (rows 14-15, copied/formula-matched
# from the former segment-7 block) plus two trailing blank-formatted rows.
# Finally, the ITEM6d sheet is made the active tab/selection (moving the
# "tabSelected" flag off ITEM6abc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITEM6d")

# --- Renumber existing segment id column (A4:A13): 3,3,4,4,5,5,6,6,7,7 -> 2,2,3,3,4,4,5,5,6,6
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 4
$ws.Range("A9").Value = 4
$ws.Range("A10").Value = 5
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 6

# --- New row 14: segment 7, "x" (same formula pattern as the other x-rows)
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "x"
$ws.Range("C14").Value = 0
$ws.Range("D14").Formula = "=2.21+0.18"
$ws.Range("E14").Formula = "=2.21+0.18+1.99"
$ws.Range("F14").Formula = "=E14+0.3+0.05"
$ws.Range("G14").Formula = "=F14+2.57+1.99"
$ws.Range("H14").Formula = "=G14+3.34+2.63"
$ws.Range("I14").Formula = "=H14+0.1+0.2"
$ws.Range("J14").Formula = "=I14+1.79+0.11+2.27"

# Copy formatting from the prior "x" row (row 12) onto the new row 14
$ws.Range("C12:J12").Copy() | Out-Null
$ws.Range("C14:J14").PasteSpecial(-4122) | Out-Null

# --- New row 15: segment 7, "z"
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "z"
$ws.Range("C15").Value = 3.29
$ws.Range("D15").Value = 1.8
$ws.Range("E15").Value = 0.2
$ws.Range("F15").Value = 0.1
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0.1
$ws.Range("I15").Value = 0.2
$ws.Range("J15").Value = 3.29

# Copy formatting from the prior "z" row (row 13) onto the new row 15
$ws.Range("C13:J13").Copy() | Out-Null
$ws.Range("C15:J15").PasteSpecial(-4122) | Out-Null

# --- Rows 16-17: trailing blank (formatted-only) rows
$ws.Range("C12:J12").Copy() | Out-Null
$ws.Range("C16:J17").PasteSpecial(-4122) | Out-Null

# --- Update the selection on ITEM6abc (previously the active tab) away from its old cell
$wsAbc = $wb.Worksheets.Item("ITEM6abc")
$wsAbc.Range("M23").Select()

# --- Make ITEM6d the active sheet/tab, with the new selection
$ws.Activate()
$ws.Range("Q15").Select()
